# Auto-generated Excel COM-interop edit script
# Applies numeric value corrections to ALC, BSM, CRP, CUL, GSM, LTW, WVR sheets
# as produced by the scheduled Sheets runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = ""

# Row 14
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = ""

# Row 17
$ws.Range("H17").Value = 1352.8846
$ws.Range("J17").Value = 1795
$ws.Range("L17").Value = 5385
$ws.Range("N17").Value = -5721

# Row 98
$ws.Range("H98").Value = 1406.6923
$ws.Range("I98").Value = 1131.7
$ws.Range("J98").Value = 2323.3333
$ws.Range("K98").Value = 1131.7
$ws.Range("L98").Value = 2323.3333
$ws.Range("M98").Value = 366.3
$ws.Range("N98").Value = -5319.3333

# Row 112
$ws.Range("H112").Value = 1903.1082
$ws.Range("I112").Value = 995
$ws.Range("K112").Value = 2985
$ws.Range("M112").Value = -1877

# Row 113
$ws.Range("H113").Value = 3754.3928
$ws.Range("I113").Value = 3036.9285
$ws.Range("J113").Value = 4471.857
$ws.Range("K113").Value = 3036.9285
$ws.Range("L113").Value = 4471.857
$ws.Range("M113").Value = 217.0715
$ws.Range("N113").Value = -10979.857

# Row 122
$ws.Range("H122").Value = 1406.6923
$ws.Range("I122").Value = 1131.7
$ws.Range("J122").Value = 2323.3333
$ws.Range("K122").Value = 3395.1
$ws.Range("L122").Value = 6969.999899999999
$ws.Range("M122").Value = -945.1000000000004
$ws.Range("N122").Value = -11869.9999

# Row 129
$ws.Range("H129").Value = 861.7727
$ws.Range("I129").Value = 417.75
$ws.Range("J129").Value = 1115.5
$ws.Range("K129").Value = 1253.25
$ws.Range("L129").Value = 3346.5
$ws.Range("M129").Value = 3746.75
$ws.Range("N129").Value = -13346.5

# Row 137
$ws.Range("H137").Value = 29155.922
$ws.Range("I137").Value = 81753.53999999999
$ws.Range("J137").Value = 1805.16
$ws.Range("K137").Value = 245260.62
$ws.Range("L137").Value = 5415.48
$ws.Range("M137").Value = -242710.62
$ws.Range("N137").Value = -10515.48

# Row 138
$ws.Range("H138").Value = 2833.255
$ws.Range("I138").Value = 5750
$ws.Range("J138").Value = 2585.0212
$ws.Range("K138").Value = 17250
$ws.Range("L138").Value = 7755.0636
$ws.Range("M138").Value = -12110
$ws.Range("N138").Value = -18035.0636

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1209.0605
$ws.Range("I94").Value = 1122.8182
$ws.Range("J94").Value = 1381.5454
$ws.Range("K94").Value = 1122.8182
$ws.Range("L94").Value = 1381.5454
$ws.Range("M94").Value = -671.8181999999999
$ws.Range("N94").Value = -2283.5454

$ws = $wb.Worksheets.Item("CRP")
# Row 19
$ws.Range("H19").Value = 50207
$ws.Range("I19").Value = 258.75
$ws.Range("J19").Value = 250000
$ws.Range("K19").Value = 258.75
$ws.Range("L19").Value = 250000
$ws.Range("M19").Value = -88.75
$ws.Range("N19").Value = -250340

# Row 24
$ws.Range("H24").Value = 50207
$ws.Range("I24").Value = 258.75
$ws.Range("J24").Value = 250000
$ws.Range("K24").Value = 258.75
$ws.Range("L24").Value = 250000
$ws.Range("M24").Value = -88.75
$ws.Range("N24").Value = -250340

# Row 31
$ws.Range("H31").Value = 2724.3088
$ws.Range("I31").Value = 1966.475
$ws.Range("J31").Value = 3806.9285
$ws.Range("K31").Value = 1966.475
$ws.Range("L31").Value = 3806.9285
$ws.Range("M31").Value = -1671.475
$ws.Range("N31").Value = -4396.9285

# Row 34
$ws.Range("H34").Value = 2724.3088
$ws.Range("I34").Value = 1966.475
$ws.Range("J34").Value = 3806.9285
$ws.Range("K34").Value = 1966.475
$ws.Range("L34").Value = 3806.9285
$ws.Range("M34").Value = -1764.475
$ws.Range("N34").Value = -4210.9285

# Row 58
$ws.Range("H58").Value = 1924.9722
$ws.Range("I58").Value = 1407.125
$ws.Range("J58").Value = 2960.6667
$ws.Range("K58").Value = 1407.125
$ws.Range("L58").Value = 2960.6667
$ws.Range("M58").Value = -1204.125
$ws.Range("N58").Value = -3366.6667

# Row 132
$ws.Range("H132").Value = 2574.9143
$ws.Range("I132").Value = 1443.3684
$ws.Range("J132").Value = 3918.625
$ws.Range("K132").Value = 4330.1052
$ws.Range("L132").Value = 11755.875
$ws.Range("M132").Value = -1800.1052
$ws.Range("N132").Value = -16815.875

# Row 136
$ws.Range("H136").Value = 1924.9722
$ws.Range("I136").Value = 1407.125
$ws.Range("J136").Value = 2960.6667
$ws.Range("K136").Value = 4221.375
$ws.Range("L136").Value = 8882.000100000001
$ws.Range("M136").Value = -1671.375
$ws.Range("N136").Value = -13982.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 546.7727
$ws.Range("I14").Value = 546.7727
$ws.Range("K14").Value = 1640.3181
$ws.Range("M14").Value = -1467.3181

# Row 38
$ws.Range("H38").Value = 49.266666
$ws.Range("I38").Value = 47.090908
$ws.Range("J38").Value = 55.25
$ws.Range("K38").Value = 141.272724
$ws.Range("L38").Value = 165.75
$ws.Range("M38").Value = 205.727276
$ws.Range("N38").Value = -859.75

# Row 92
$ws.Range("H92").Value = 879
$ws.Range("I92").Value = 840
$ws.Range("J92").Value = 918
$ws.Range("K92").Value = 2520
$ws.Range("L92").Value = 2754
$ws.Range("M92").Value = -1272
$ws.Range("N92").Value = -5250

# Row 125
$ws.Range("H125").Value = 2847
$ws.Range("I125").Value = 430
$ws.Range("J125").Value = 3249.8333
$ws.Range("K125").Value = 1290
$ws.Range("L125").Value = 9749.499899999999
$ws.Range("M125").Value = 3630
$ws.Range("N125").Value = -19589.4999

# Row 131
$ws.Range("H131").Value = 864.13336
$ws.Range("I131").Value = 556.8
$ws.Range("J131").Value = 911.4154
$ws.Range("K131").Value = 1670.4
$ws.Range("L131").Value = 2734.2462
$ws.Range("M131").Value = 3369.6
$ws.Range("N131").Value = -12814.2462

$ws = $wb.Worksheets.Item("GSM")
# Row 13
$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = -1278

# Row 122
$ws.Range("H122").Value = 1831
$ws.Range("I122").Value = 1730.1428
$ws.Range("J122").Value = 2066.3333
$ws.Range("K122").Value = 5190.428400000001
$ws.Range("L122").Value = 6198.999899999999
$ws.Range("M122").Value = -2740.428400000001
$ws.Range("N122").Value = -11098.9999

# Row 132
$ws.Range("H132").Value = 3219.9822
$ws.Range("I132").Value = 3089.6667
$ws.Range("J132").Value = 3454.55
$ws.Range("K132").Value = 9269.000100000001
$ws.Range("L132").Value = 10363.65
$ws.Range("M132").Value = -6739.000100000001
$ws.Range("N132").Value = -15423.65

$ws = $wb.Worksheets.Item("LTW")
# Row 19
$ws.Range("H19").Value = 790
$ws.Range("I19").Value = 500
$ws.Range("J19").Value = 935
$ws.Range("K19").Value = 500
$ws.Range("L19").Value = 935
$ws.Range("M19").Value = -330
$ws.Range("N19").Value = -1275

# Row 40
$ws.Range("H40").Value = 4027.6365
$ws.Range("I40").Value = 4033.7778
$ws.Range("K40").Value = 4033.7778
$ws.Range("M40").Value = -3897.7778

# Row 46
$ws.Range("H46").Value = 766.2069
$ws.Range("I46").Value = 705.2174
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 705.2174
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -517.2174
$ws.Range("N46").Value = -1376

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 52011.95
$ws.Range("I122").Value = 73189.92999999999
$ws.Range("J122").Value = 2596.6667
$ws.Range("K122").Value = 219569.79
$ws.Range("L122").Value = 7790.000100000001
$ws.Range("M122").Value = -217119.79
$ws.Range("N122").Value = -12690.0001

# Row 126
$ws.Range("H126").Value = 1565.375
$ws.Range("I126").Value = 1192.2727
$ws.Range("J126").Value = 2386.2
$ws.Range("K126").Value = 3576.8181
$ws.Range("L126").Value = 7158.599999999999
$ws.Range("M126").Value = -1106.8181
$ws.Range("N126").Value = -12098.6
